$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.821390151977539
$ws.Range("B1").Value = 2.770732879638672
$ws.Range("C1").Value = 2.998645067214966
$ws.Range("D1").Value = 3.523798942565918
$ws.Range("E1").Value = 3.204209327697754
